$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 105, pushing existing rows 105-165 down to 106-166
$ws.Rows.Item(105).Insert()

$ws.Cells.Item(105, 1).Value = 4
$ws.Cells.Item(105, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(105, 3).Value = "Los Lagos"
$ws.Cells.Item(105, 4).Value = 45202
$ws.Cells.Item(105, 5).Value = 10
$ws.Cells.Item(105, 6).Value = 100112026
$ws.Cells.Item(105, 7).Value = "Haba"
$ws.Cells.Item(105, 8).Value = "Sin especificar"
$ws.Cells.Item(105, 9).Value = "Primera"
$ws.Cells.Item(105, 10).Value = 80
$ws.Cells.Item(105, 11).Value = 15000
$ws.Cells.Item(105, 12).Value = 15000
$ws.Cells.Item(105, 13).Value = 15000
$ws.Cells.Item(105, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(105, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(105, 16).Value = 600
$ws.Cells.Item(105, 17).Value = 25
$ws.Cells.Item(105, 18).Value = "Hortaliza"
